# Applies the "Optuna Attempt (go back with original)" edit:
# Updates Seasonality Index (L) and MyForecast (D) values on the
# "Forecast Comparison" sheet, and the dependent summary figures on
# the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Seasonality Index (column L) updates
$wsForecast.Range("L2").Value  = 0.82
$wsForecast.Range("L3").Value  = 0.85
$wsForecast.Range("L4").Value  = 1.12
$wsForecast.Range("L5").Value  = 1.04
$wsForecast.Range("L6").Value  = 0.83
$wsForecast.Range("L7").Value  = 1.09
$wsForecast.Range("L8").Value  = 0.89
$wsForecast.Range("L9").Value  = 0.82
$wsForecast.Range("L10").Value = 1.15
$wsForecast.Range("L11").Value = 1.07
$wsForecast.Range("L12").Value = 1.16
$wsForecast.Range("L13").Value = 0.92
$wsForecast.Range("L14").Value = 1
$wsForecast.Range("L15").Value = 0.82
$wsForecast.Range("L16").Value = 1.16
$wsForecast.Range("L17").Value = 0.88

# MyForecast (column D) updates
$wsForecast.Range("D13").Value = 21
$wsForecast.Range("D16").Value = 19
$wsForecast.Range("D17").Value = 17

# Summary sheet updates (stored as text values, matching source data type)
$wsSummary.Range("B9").NumberFormat  = "@"
$wsSummary.Range("B9").Value         = "342"
$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value        = "17"
